$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1, col A)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 06:22"

# Alemania (row 8)
$ws.Range("D8").Value = 109800
$ws.Range("E8").Value = 39439

# Brasil (row 14)
$ws.Range("B14").Value = 54043
$ws.Range("C14").Value = 1048
$ws.Range("E14").Value = 22684
$ws.Range("G14").Value = 34
$ws.Range("H14").Value = 3704

# Suecia (row 24)
$ws.Range("D24").Value = 1005
$ws.Range("E24").Value = 14410

# Mexico (row 28)
$ws.Range("D28").Value = 7149
$ws.Range("E28").Value = 4502

# Japon (row 29)
$ws.Range("B29").Value = 12829
$ws.Range("C29").Value = 117
$ws.Range("E29").Value = 10954

# Australia (row 46)
$ws.Range("B46").Value = 6692
$ws.Range("C46").Value = 17
$ws.Range("D46").Value = 5357
$ws.Range("E46").Value = 1255
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 80

# Panama (row 49)
$ws.Range("B49").Value = 5338
$ws.Range("C49").Value = 172
$ws.Range("D49").Value = 319
$ws.Range("E49").Value = 4865
$ws.Range("F49").Value = 87
$ws.Range("G49").Value = 8
$ws.Range("H49").Value = 154

# Swap Honduras/Uruguay order: row 102 becomes Honduras, row 103 becomes Uruguay,
# and the province/country figures are refreshed for both.
$ws.Range("A102").Value = "Honduras"
$ws.Range("B102").Value = 591
$ws.Range("C102").Value = 29
$ws.Range("D102").Value = 58
$ws.Range("E102").Value = 478
$ws.Range("F102").Value = 10
$ws.Range("G102").Value = 8
$ws.Range("H102").Value = 55

$ws.Range("A103").Value = "Uruguay"
$ws.Range("B103").Value = 563
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 369
$ws.Range("E103").Value = 182
$ws.Range("F103").Value = 9
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 12

# Nepal (row 166)
$ws.Range("D166").Value = 11
$ws.Range("E166").Value = 38

# Mongolia (row 173)
$ws.Range("B173").Value = 37
$ws.Range("C173").Value = 1
$ws.Range("E173").Value = 28
